# Rename the worksheet from "testGoogle" to "Sheet1"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet1"

# Move the active selection from B5 to C4
$ws.Range("C4").Select() | Out-Null

# Swap the theme's Accent1 / Accent5 colors (5B9BD5 <-> 4472C4)
$theme = $wb.Theme
$colorScheme = $theme.ThemeColorScheme
$colorScheme.Colors(5).RGB = 12874308   # msoThemeAccent1 -> 4472C4
$colorScheme.Colors(9).RGB = 13998939   # msoThemeAccent5 -> 5B9BD5
